$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.100.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.31%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.909.56'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.41%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.69%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.08%  '

$ws.Range("E7").Value = '  +3.56%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.89'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.44%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2971'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.60%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06823'
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.911.46'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.55%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '17.26'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.39%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07364'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.56%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6926'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.65%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '86.72'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.87%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.873'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.45%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.107.88'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.42%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008226'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +11.04%  '

$ws.Range("E19").Value = '  +0.06%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.99'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.92%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.155.99'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.28%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.815'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.61%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.716'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.84%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.189'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.50%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '146.79'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '135.24'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.19%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.64%  '

$ws.Range("E29").Value = '  +5.35%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.392'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.00%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.232'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.50%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08819'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.37%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.009'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.46%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05065'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.67%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.145'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.70%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7157'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.29%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.690'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.48%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.811'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.82%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.269'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.51%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9657'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.78%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01692'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.02%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.161'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.94%  '

$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '105.39'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.41%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4298'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9990'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.02%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.639'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.63%  '

$ws.Range("E47").Value = '  +4.23%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05735'
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '33.12'
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.416'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.49%  '

$ws.Range("E51").Value = '  +4.48%  '

Write-Host "Update complete"